$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each functionality row, the "Attributes" column (D) held a string of the
# form ",[before...],[after...]" (a stray leading comma, and a comma glueing
# the two attribute lists together). Separate the data elements with a comma
# for each functionality instead: drop the leading comma and concatenate the
# two attribute-list reprs directly; when a cell bundles more than one
# functionality's data, keep a single comma between functionality chunks.

$ws.Range("D2").Value  = "['Amt_avail', 'From_AcctNum', 'To_AcctNum', 'Acc_type']['Amt_avail', 'Acc_type']"
$ws.Range("D3").Value  = "['Bill_type', 'Max_limit']['Amt_avail', 'From_AcctNum', 'To_AcctNum', 'Acc_type']"
$ws.Range("D4").Value  = "['From_AcctNum', 'To_AcctNum', 'Amt_avail', 'Amt_trnsfr']['Amt_avail', 'From_AcctNum', 'To_AcctNum', 'Acc_type']"
$ws.Range("D5").Value  = "['Cus_Nme', 'Amt_avail', 'Acc_type']['Bill_type', 'Max_limit']"
$ws.Range("D6").Value  = "['Loan_Amt', 'Amt_avail', 'Debit_pin', 'Amt_wdrl', 'Amt_trnsfr']['Loan_Amt', 'Cred_Score']"
$ws.Range("D7").Value  = "[]['Amt_avail', 'From_AcctNum', 'To_AcctNum', 'Acc_type']"
$ws.Range("D8").Value  = "['Cus_Nme', 'Acc_type']['Bill_type', 'Max_limit']"
$ws.Range("D9").Value  = "['Cus_Nme', 'Acc_type']['Bill_type', 'Max_limit']"
$ws.Range("D10").Value = "[]['Amt_avail']"
$ws.Range("D11").Value = "['Loan_Amt', 'Cred_Score'][],['Loan_Amt', 'Amt_avail', 'Debit_pin', 'Amt_wdrl', 'Amt_trnsfr']['Loan_Amt', 'Cred_Score']"

# D12 never had real attribute data (empty placeholder cell) - clear it out.
$ws.Range("D12").Value = ""
